$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, "H").Value = 522.625
$ws.Cells.Item(28, "I").Value = 422.07693
$ws.Cells.Item(28, "J").Value = 958.3333
$ws.Cells.Item(28, "K").Value = 422.07693
$ws.Cells.Item(28, "L").Value = 958.3333
$ws.Cells.Item(28, "M").Value = 62.92307
$ws.Cells.Item(28, "N").Value = -1928.3333
$ws.Cells.Item(40, "H").Value = 5209.385
$ws.Cells.Item(40, "I").Value = 2744.4
$ws.Cells.Item(40, "K").Value = 2744.4
$ws.Cells.Item(40, "M").Value = -2569.4
$ws.Cells.Item(70, "H").Value = 3095.5557
$ws.Cells.Item(70, "I").Value = 2190.3333
$ws.Cells.Item(70, "K").Value = 6570.999899999999
$ws.Cells.Item(70, "M").Value = -6300.999899999999
$ws.Cells.Item(73, "H").Value = 3095.5557
$ws.Cells.Item(73, "I").Value = 2190.3333
$ws.Cells.Item(73, "K").Value = 6570.999899999999
$ws.Cells.Item(73, "M").Value = -5634.999899999999
$ws.Cells.Item(76, "H").Value = 8074.5
$ws.Cells.Item(76, "I").Value = 900
$ws.Cells.Item(76, "K").Value = 900
$ws.Cells.Item(76, "M").Value = -585
$ws.Cells.Item(79, "H").Value = 8074.5
$ws.Cells.Item(79, "I").Value = 900
$ws.Cells.Item(79, "K").Value = 900
$ws.Cells.Item(79, "M").Value = 192
$ws.Cells.Item(116, "H").Value = 3872.6365
$ws.Cells.Item(116, "I").Value = 3474.875
$ws.Cells.Item(116, "K").Value = 3474.875
$ws.Cells.Item(116, "M").Value = -32.875
$ws.Cells.Item(131, "H").Value = 148373.42
$ws.Cells.Item(131, "I").Value = 202577.8
$ws.Cells.Item(131, "K").Value = 607733.3999999999
$ws.Cells.Item(131, "M").Value = -602693.3999999999
$ws.Cells.Item(132, "H").Value = 7472
$ws.Cells.Item(132, "I").Value = 8296.214
$ws.Cells.Item(132, "K").Value = 24888.642
$ws.Cells.Item(132, "M").Value = -22358.642
$ws.Cells.Item(137, "H").Value = 2038.4445
$ws.Cells.Item(137, "I").Value = 1558.5
$ws.Cells.Item(137, "J").Value = 2998.3333
$ws.Cells.Item(137, "K").Value = 4675.5
$ws.Cells.Item(137, "L").Value = 8994.999899999999
$ws.Cells.Item(137, "M").Value = -2125.5
$ws.Cells.Item(137, "N").Value = -14094.9999
$ws.Cells.Item(138, "H").Value = 2591.7678
$ws.Cells.Item(138, "J").Value = 3015.8108
$ws.Cells.Item(138, "L").Value = 9047.432400000002
$ws.Cells.Item(138, "N").Value = -19327.4324
$ws.Cells.Item(141, "H").Value = 6652.467
$ws.Cells.Item(141, "I").Value = 5819.8335
$ws.Cells.Item(141, "K").Value = 17459.5005
$ws.Cells.Item(141, "M").Value = -12279.5005

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, "H").Value = 4222.857
$ws.Cells.Item(122, "I").Value = 4222.857
$ws.Cells.Item(122, "K").Value = 12668.571
$ws.Cells.Item(122, "M").Value = -10218.571
$ws.Cells.Item(138, "H").Value = 199999
$ws.Cells.Item(138, "J").Value = 199999
$ws.Cells.Item(138, "L").Value = 199999
$ws.Cells.Item(138, "N").Value = -210279

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, "H").Value = 3532.8298
$ws.Cells.Item(134, "I").Value = 3532.8298
$ws.Cells.Item(134, "K").Value = 10598.4894
$ws.Cells.Item(134, "M").Value = -8063.4894

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(6, "H").Value = 1486
$ws.Cells.Item(6, "J").Value = 1486
$ws.Cells.Item(6, "L").Value = 1486
$ws.Cells.Item(6, "N").Value = -1712
$ws.Cells.Item(31, "H").Value = 4766.567
$ws.Cells.Item(31, "I").Value = 3337.6365
$ws.Cells.Item(31, "J").Value = 8696.125
$ws.Cells.Item(31, "K").Value = 3337.6365
$ws.Cells.Item(31, "L").Value = 8696.125
$ws.Cells.Item(31, "M").Value = -3042.6365
$ws.Cells.Item(31, "N").Value = -9286.125
$ws.Cells.Item(34, "H").Value = 4766.567
$ws.Cells.Item(34, "I").Value = 3337.6365
$ws.Cells.Item(34, "J").Value = 8696.125
$ws.Cells.Item(34, "K").Value = 3337.6365
$ws.Cells.Item(34, "L").Value = 8696.125
$ws.Cells.Item(34, "M").Value = -3135.6365
$ws.Cells.Item(34, "N").Value = -9100.125
$ws.Cells.Item(74, "H").Value = 41870
$ws.Cells.Item(74, "J").Value = 41870
$ws.Cells.Item(74, "L").Value = 41870
$ws.Cells.Item(74, "N").Value = -43618
$ws.Cells.Item(77, "H").Value = 41870
$ws.Cells.Item(77, "J").Value = 41870
$ws.Cells.Item(77, "L").Value = 125610
$ws.Cells.Item(77, "N").Value = -134346
$ws.Cells.Item(134, "H").Value = 4796.524
$ws.Cells.Item(134, "I").Value = 2997.3635
$ws.Cells.Item(134, "J").Value = 6775.6
$ws.Cells.Item(134, "K").Value = 8992.0905
$ws.Cells.Item(134, "L").Value = 20326.8
$ws.Cells.Item(134, "M").Value = -6457.0905
$ws.Cells.Item(134, "N").Value = -25396.8
$ws.Cells.Item(141, "H").Value = 37850
$ws.Cells.Item(141, "J").Value = 37850
$ws.Cells.Item(141, "L").Value = 37850
$ws.Cells.Item(141, "N").Value = -48210

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, "H").Value = 9601354
$ws.Cells.Item(4, "I").Value = 3281224.5
$ws.Cells.Item(4, "K").Value = 9843673.5
$ws.Cells.Item(4, "M").Value = -9843561.5
$ws.Cells.Item(5, "H").Value = 1228.1111
$ws.Cells.Item(5, "I").Value = 208.83333
$ws.Cells.Item(5, "J").Value = 3266.6667
$ws.Cells.Item(5, "K").Value = 626.49999
$ws.Cells.Item(5, "L").Value = 9800.000100000001
$ws.Cells.Item(5, "M").Value = -514.49999
$ws.Cells.Item(5, "N").Value = -10024.0001
$ws.Cells.Item(8, "H").Value = 725
$ws.Cells.Item(8, "I").Value = 725
$ws.Cells.Item(8, "K").Value = 2175
$ws.Cells.Item(8, "M").Value = -2036
$ws.Cells.Item(11, "H").Value = 10000650
$ws.Cells.Item(11, "I").Value = 14286214
$ws.Cells.Item(11, "K").Value = 42858642
$ws.Cells.Item(11, "M").Value = -42858502
$ws.Cells.Item(86, "H").Value = 420.77777
$ws.Cells.Item(86, "I").Value = 422
$ws.Cells.Item(86, "J").Value = 420.42856
$ws.Cells.Item(86, "K").Value = 1266
$ws.Cells.Item(86, "L").Value = 1261.28568
$ws.Cells.Item(86, "M").Value = -80
$ws.Cells.Item(86, "N").Value = -3633.28568
$ws.Cells.Item(89, "H").Value = 420.77777
$ws.Cells.Item(89, "I").Value = 422
$ws.Cells.Item(89, "J").Value = 420.42856
$ws.Cells.Item(89, "K").Value = 3798
$ws.Cells.Item(89, "L").Value = 3783.85704
$ws.Cells.Item(89, "M").Value = 2130
$ws.Cells.Item(89, "N").Value = -15639.85704
$ws.Cells.Item(131, "H").Value = 3705341.2
$ws.Cells.Item(131, "I").Value = 499.76923
$ws.Cells.Item(131, "K").Value = 1499.30769
$ws.Cells.Item(131, "M").Value = 3540.69231
$ws.Cells.Item(135, "H").Value = 1228.1111
$ws.Cells.Item(135, "I").Value = 208.83333
$ws.Cells.Item(135, "J").Value = 3266.6667
$ws.Cells.Item(135, "K").Value = 1879.49997
$ws.Cells.Item(135, "L").Value = 29400.0003
$ws.Cells.Item(135, "M").Value = 655.5000300000002
$ws.Cells.Item(135, "N").Value = -34470.0003

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, "H").Value = 3878.2693
$ws.Cells.Item(2, "I").Value = 22.473684
$ws.Cells.Item(2, "K").Value = 22.473684
$ws.Cells.Item(2, "M").Value = 90.52631600000001
$ws.Cells.Item(80, "H").Value = 4247.2856
$ws.Cells.Item(80, "I").Value = 3990.8
$ws.Cells.Item(80, "J").Value = 4888.5
$ws.Cells.Item(80, "K").Value = 3990.8
$ws.Cells.Item(80, "L").Value = 4888.5
$ws.Cells.Item(80, "M").Value = -2992.8
$ws.Cells.Item(80, "N").Value = -6884.5
$ws.Cells.Item(83, "H").Value = 4247.2856
$ws.Cells.Item(83, "I").Value = 3990.8
$ws.Cells.Item(83, "J").Value = 4888.5
$ws.Cells.Item(83, "K").Value = 19954
$ws.Cells.Item(83, "L").Value = 24442.5
$ws.Cells.Item(83, "M").Value = -14962
$ws.Cells.Item(83, "N").Value = -34426.5
$ws.Cells.Item(122, "H").Value = 2562.9092
$ws.Cells.Item(122, "I").Value = 2519.3
$ws.Cells.Item(122, "J").Value = 2999
$ws.Cells.Item(122, "K").Value = 7557.900000000001
$ws.Cells.Item(122, "L").Value = 8997
$ws.Cells.Item(122, "M").Value = -5107.900000000001
$ws.Cells.Item(122, "N").Value = -13897
$ws.Cells.Item(135, "H").Value = 123750
$ws.Cells.Item(135, "J").Value = 123750
$ws.Cells.Item(135, "L").Value = 123750
$ws.Cells.Item(135, "N").Value = -133890

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, "H").Value = 624.625
$ws.Cells.Item(55, "I").Value = 722.2143
$ws.Cells.Item(55, "K").Value = 722.2143
$ws.Cells.Item(55, "M").Value = -549.2143
$ws.Cells.Item(93, "H").Value = 17036.846
$ws.Cells.Item(93, "I").Value = 1043.6364
$ws.Cells.Item(93, "K").Value = 1043.6364
$ws.Cells.Item(93, "M").Value = 204.3635999999999
$ws.Cells.Item(136, "H").Value = 4709.2915
$ws.Cells.Item(136, "I").Value = 4556.8887
$ws.Cells.Item(136, "K").Value = 13670.6661
$ws.Cells.Item(136, "M").Value = -11120.6661

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, "H").Value = 9999.333000000001
$ws.Cells.Item(62, "J").Value = 9999.333000000001
$ws.Cells.Item(62, "L").Value = 9999.333000000001
$ws.Cells.Item(62, "N").Value = -11247.333
$ws.Cells.Item(65, "H").Value = 9999.333000000001
$ws.Cells.Item(65, "J").Value = 9999.333000000001
$ws.Cells.Item(65, "L").Value = 49996.665
$ws.Cells.Item(65, "N").Value = -56236.665
$ws.Cells.Item(107, "H").Value = 568.375
$ws.Cells.Item(107, "I").Value = 544.9091
$ws.Cells.Item(107, "K").Value = 1634.7273
$ws.Cells.Item(107, "M").Value = 285.2727
$ws.Cells.Item(136, "H").Value = 4481.4346
$ws.Cells.Item(136, "I").Value = 3884.9697
$ws.Cells.Item(136, "K").Value = 11654.9091
$ws.Cells.Item(136, "M").Value = -9104.909100000001
